# north-virginia-2024.xlsx: "Changed up the counties."
# Update the per-county "result" values (column C) on the only worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = -57.9
$ws.Range("C12").Value = -68.2
$ws.Range("C20").Value = -57.8
$ws.Range("C24").Value = -54.9
$ws.Range("C32").Value = 54.8
$ws.Range("C35").Value = 54.7
$ws.Range("C44").Value = -57.9
$ws.Range("C48").Value = 52.7
$ws.Range("C51").Value = -53.7
$ws.Range("C54").Value = -74.3

# Reset the view: scroll back to the top (drops the saved topLeftCell="A46")
# and move the active selection to V17, matching the new saved view state.
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("V17").Select()
